# "Generate Report for Handback" — record the handback completion for the
# zh-cn and de-de localization targets: the overview status flips from
# "Ready for handoff" to "Handed back: in sync with en-US", each language
# sheet gets its generated-target markdown linked (with the hyperlink style)
# and its handback .xlf filename recorded, and the de-de sheet's handback
# timestamp is refreshed. A few columns also widen to fit the new long
# filename/link text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: status text + wider zh-cn/de-de status columns
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"

$ov.Range("E1").ColumnWidth = 29.17
$ov.Range("F1").ColumnWidth = 29.17

# ---------------------------------------------------------------------
# zh-cn sheet: link the generated target file, record the handback xlf
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C1").ColumnWidth = 29.17
$zh.Range("I1").ColumnWidth = 39.17
$zh.Range("J1").ColumnWidth = 39.17

$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71ba8a9f0c3abc1f7e2064806fda52c4385dfe89/e2e/a2da014b-a8f9-4336-9e36-51b3caa38b79.md", "", "", "a2da014b-a8f9-4336-9e36-51b3caa38b79.md") | Out-Null
$zh.Range("I2").Style = "HyperLink"

$zh.Range("J2").Value = "a2da014b-a8f9-4336-9e36-51b3caa38b79.655a220044a1cbe7945451f23240f66a31b5d86e.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-30 09:23:46"

# ---------------------------------------------------------------------
# de-de sheet: link the generated target file, record the handback xlf
# and refresh the handback timestamp
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C1").ColumnWidth = 29.17
$de.Range("I1").ColumnWidth = 39.17
$de.Range("J1").ColumnWidth = 39.17

$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71ba8a9f0c3abc1f7e2064806fda52c4385dfe89/e2e/a2da014b-a8f9-4336-9e36-51b3caa38b79.md", "", "", "a2da014b-a8f9-4336-9e36-51b3caa38b79.md") | Out-Null
$de.Range("I2").Style = "HyperLink"

$de.Range("J2").Value = "a2da014b-a8f9-4336-9e36-51b3caa38b79.655a220044a1cbe7945451f23240f66a31b5d86e.de-de.xlf"
$de.Range("K2").Value = "2016-08-30 09:24:11"
